$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J16").Value = 35.58
$ws.Range("J17").Value = 39.76
$ws.Range("J18").Value = 37.54
$ws.Range("J19").Value = 54.1

$ws.Range("M18").Select()
